$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(104).Insert()

$ws.Cells.Item(104, 1).Value = 7
$ws.Cells.Item(104, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(104, 3).Value = "Ñuble"
$ws.Cells.Item(104, 4).Value = 44985
$ws.Cells.Item(104, 5).Value = 16
$ws.Cells.Item(104, 6).Value = 100112021
$ws.Cells.Item(104, 7).Value = "Ají"
$ws.Cells.Item(104, 8).Value = "Cristal"
$ws.Cells.Item(104, 9).Value = "Primera"
$ws.Cells.Item(104, 10).Value = 30
$ws.Cells.Item(104, 11).Value = 15000
$ws.Cells.Item(104, 12).Value = 15000
$ws.Cells.Item(104, 13).Value = 15000
$ws.Cells.Item(104, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(104, 15).Value = "Región del Maule"
$ws.Cells.Item(104, 16).Value = 600
$ws.Cells.Item(104, 17).Value = 25
$ws.Cells.Item(104, 18).Value = "Hortaliza"
